$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Fitness column (C2:C252) to the corrected value of 7569
$ws.Range("C2:C252").Value = 7569
